$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values scraped from the updated simulation run ("case with 380 kV done").
# Only the loading_percent values in columns B, D, E, F, G, L, M, O for rows 2-25 changed;
# column A (index) and the always-zero columns (C, H, I, J, K, N) are untouched.
$data = @{}
$data["B"] = @(18.29285891475653, 17.91002781494124, 17.67178514353304, 17.57401109913511, 17.55773758621723, 17.67046916398255, 18.16158525333218, 19.09471769727332, 19.75612107550752, 20.05062356980073, 20.16114049524623, 20.1373847263835, 20.05973648497879, 20.01204132740975, 19.73673870933031, 19.56615234578071, 19.46744059489957, 19.43391909455731, 19.58437379346802, 20.08257159835923, 20.40227942883905, 20.23221215987532, 19.57613786161398, 18.84611854726164)
$data["D"] = @(3.711001506839609, 3.675709904258966, 3.654094645468905, 3.645306830164832, 3.643849067622832, 3.653976037367136, 3.698825982275109, 3.786896867487498, 3.851266720933376, 3.880391730636763, 3.891390958897044, 3.88902351558862, 3.881297281237023, 3.87656065058788, 3.849359553492927, 3.832627155740714, 3.822988751922844, 3.819723106856336, 3.834409887862217, 3.883567528306816, 3.915518078645882, 3.898483980648014, 3.833603972413276, 3.76310093047399)
$data["E"] = @(31.41613300032321, 30.37976779512441, 29.72890984911355, 29.46044376186858, 29.41568259000698, 29.72530174805833, 31.06204766964069, 33.55169034791422, 35.28101464884764, 36.04261700474061, 36.32716652254089, 36.26605866503191, 36.06610558971624, 35.94311998938817, 35.23071576955347, 34.78707114361333, 34.52955101278329, 34.44196395114501, 34.83454270864097, 36.12494303976107, 36.94574311718196, 36.50980417540526, 34.81308846661614, 32.89442825646503)
$data["F"] = @(19.87892122332325, 19.89391527156841, 19.9142119897617, 19.92525137699649, 19.92725099180656, 19.91434969172778, 19.88177818131426, 19.90663366700422, 19.97971164775285, 20.02492859268647, 20.04377197821304, 20.03963725116539, 20.02644438499297, 20.01858735651717, 19.97699755745581, 19.95455085058497, 19.9427674272231, 19.93897134430867, 19.95682366050356, 20.03027277959296, 20.08830213858106, 20.05641483831119, 19.95579262965819, 19.89031989593194)
$data["G"] = @(3.592876019060736, 3.59552839431446, 3.597240754916152, 3.597959698021927, 3.598080357037946, 3.597250365134844, 3.5937732107276, 3.587616060943575, 3.583491045095374, 3.581700035747637, 3.581034043394305, 3.581176934177416, 3.581644999557491, 3.58193329303047, 3.583609806041425, 3.584660137701982, 3.585272310620566, 3.585480966514668, 3.584547495459156, 3.58150718619388, 3.579591387889319, 3.580607391314766, 3.584598395084276, 3.589211392549477)
$data["L"] = @(10.63382016981005, 10.30747130756335, 10.10108660454627, 10.01556898918037, 10.00128640852304, 10.09993887717696, 10.52259638158697, 11.29999107936066, 11.8352921558963, 12.07023511799338, 12.15791304708702, 12.13908825340145, 12.07747459233991, 12.03956486512339, 11.81976057505326, 11.68267999311321, 11.60303140106847, 11.57592763082783, 11.69735605059014, 12.09560745249786, 12.34834717894843, 12.21416259284427, 11.69072361430616, 11.0956864988069)
$data["M"] = @(16.77267393523984, 16.48511143800179, 16.30622362956135, 16.23281112785178, 16.22059205836518, 16.30523555140169, 16.67404591697704, 17.37588573200204, 17.87469431156004, 18.09721313902012, 18.18078533453217, 18.16281820268772, 18.10410284779788, 18.06804632284555, 17.86005866402114, 17.73129690307066, 17.65682652071764, 17.63154356836083, 17.74504668832571, 18.12136817700009, 18.36326400050141, 18.23454915793447, 17.73883178731289, 17.1887118311625)
$data["O"] = @(17.59287449200634, 17.68772816083751, 17.75462494755867, 17.784040142523, 17.7890539527958, 17.75501295911342, 17.62377012219742, 17.43605606307634, 17.34187175548367, 17.30877134085245, 17.29765727701068, 17.29998742878302, 17.30782840944815, 17.31281675245303, 17.34423269081632, 17.3660151338374, 17.37945995055078, 17.38416887873576, 17.3636013969979, 17.3054866297196, 17.27579206417633, 17.29087639530636, 17.3646897774189, 17.47924651002011)

$cols = @("B", "D", "E", "F", "G", "L", "M", "O")

foreach ($col in $cols) {
    $colValues = $data[$col]
    for ($i = 0; $i -lt $colValues.Length; $i++) {
        $row = $i + 2
        $ws.Range("$col$row").Value = $colValues[$i]
    }
}
